# Weekly driver report update for 2025-04-28
# Updates the "Driver Summary" sheet:
#  - Bad Drivers table: refresh counts for the remaining two drivers and
#    drop the row for "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2" (its
#    numbers move into the row above it), then refresh the Totals row.
#  - Good Drivers table: drop the "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
#    row, refresh sample counts for the two rows that shift up, and drop all
#    of the old AX201 driver rows at the bottom (only the three AX211 rows
#    remain).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -----------------------------------------------

# Row 3 ("Fi - 16.0 (1657)") gets refreshed Critical Minutes / Good Roaming values.
$ws.Range("C3").Value = 162
$ws.Range("D3").Value = 95.3

# Row 4 ("Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2") is removed entirely;
# this shifts the AX211 23.90.0.2 row (old row 5) and the Totals row
# (old row 6) up by one.
$ws.Rows(4).Delete()

# The AX211 160MHz - 23.90.0.2 row (now row 4) gets refreshed numbers.
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 96.90000000000001

# The Totals row (now row 5) gets refreshed Client Count / Critical Minutes.
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 168

# --- Good Drivers table ------------------------------------------------

# "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4" (row 13 after the shift above)
# is removed; the three AX211 rows below it shift up by one.
$ws.Rows(13).Delete()

# Refresh Total Samples for the two rows that shifted (22.150.3.1 / 22.150.0.3).
$ws.Range("B13").Value = 11140
$ws.Range("B14").Value = 14487

# The remaining five legacy AX201 rows are dropped entirely.
$ws.Rows("16:20").Delete()
